$d = $word.ActiveDocument

function Find-TargetIndex {
    # Locate the first existing checklist item ("Código completado") - the
    # two new bullet points are inserted right before it, i.e. directly
    # after the introductory paragraph ("Un story en estado Verify...").
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*digo completado*") {
            return $i
        }
    }
    return -1
}

# --- Insert paragraph #1: "Story comenzada" ---
$idx = Find-TargetIndex
$r = $d.Paragraphs($idx).Range
$r.Collapse(1)   # wdCollapseStart
$r.InsertParagraphBefore()

$idx = Find-TargetIndex
$newPara1 = $d.Paragraphs($idx - 1)
$newPara1.Range.Text = "Story comenzada"

# --- Insert paragraph #2: "No existencia de bugs de importancia crítica o alta" ---
$idx = Find-TargetIndex
$r = $d.Paragraphs($idx).Range
$r.Collapse(1)
$r.InsertParagraphBefore()

$idx = Find-TargetIndex
$newPara2 = $d.Paragraphs($idx - 1)

# The paragraph also carries Word's "_GoBack" bookmark (last-edit-location
# marker) right after its text. A range collapsed exactly at a paragraph's
# end sits on an ambiguous boundary, so a trailing sentinel character is
# typed first to give the bookmark an unambiguous, interior anchor point;
# the sentinel is removed again immediately after.
$realText = "No existencia de bugs de importancia crítica o alta"
$sentinel = "#"
$newPara2.Range.Text = $realText + $sentinel

$idx = Find-TargetIndex
$newPara2 = $d.Paragraphs($idx - 1)
$bmPos = $newPara2.Range.Start + $realText.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the sentinel character again.
$sentinelRange = $d.Range($bmPos, $bmPos + 1)
$sentinelRange.Text = ""
